$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "release/8.0.16"
$ws.Range("B19").Value = "X"
$ws.Range("C19").Value = "X"
$ws.Range("D19").Value = "X"
$ws.Range("E19").Value = "X"
